# Ratios_add.xlsx: add a "Lab. #" column at the front of the table and
# highlight the data row with a light-green fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before column A; this shifts all the existing
# "Ratio ..." / "Error (%) ..." columns one place to the right.
$ws.Columns("A:A").Insert()

# New header + laboratory number for the data row.
$ws.Range("A1").Value = "Lab. #"
$ws.Range("A2").Value = 8166

# The lab-number column is much narrower than the ratio columns.
$ws.Columns("A:A").ColumnWidth = 6.9

# Highlight the whole data row (now A2:Y2) with the light olive-green fill
# used to flag corrected rows.
$ws.Range("A2:Y2").Interior.Color = 12379352

Write-Output "done"
